$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.444.99"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "'1.655.14"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'307.99"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'0.9987"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.3620"
$ws.Range("E7").Value = "  -4.01%  "
$ws.Range("D8").Value = "'47.33"
$ws.Range("E8").Value = "  -4.40%  "
$ws.Range("D9").Value = "'0.3280"
$ws.Range("E9").Value = "  -5.62%  "
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").Value = "'0.06949"
$ws.Range("E11").Value = "  -6.98%  "
$ws.Range("D12").Value = "'0.9985"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "'5.959"
$ws.Range("E13").Value = "  -5.18%  "
$ws.Range("D14").Value = "'19.33"
$ws.Range("E14").Value = "  -7.43%  "
$ws.Range("D15").Value = "'6.629"
$ws.Range("E15").Value = "  -4.82%  "
$ws.Range("D16").Value = "'1.650.79"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "'0.00001042"
$ws.Range("E17").Value = "  -7.17%  "
$ws.Range("D18").Value = "'0.06520"
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "'0.9978"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'76.56"
$ws.Range("E20").Value = "  -8.90%  "
$ws.Range("D21").Value = "'5.922"
$ws.Range("E21").Value = "  -7.18%  "
$ws.Range("D22").Value = "'15.73"
$ws.Range("E22").Value = "  -8.56%  "
$ws.Range("D23").Value = "'12.59"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "'24.401.45"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'2.429"
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "'2.348"
$ws.Range("E26").Value = "  -15.68%  "
$ws.Range("D27").Value = "'146.61"
$ws.Range("E27").Value = "  -2.85%  "
$ws.Range("D28").Value = "'18.36"
$ws.Range("E28").Value = "  -10.32%  "
$ws.Range("D29").Value = "'1.840.42"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").Value = "'124.32"
$ws.Range("E30").Value = "  -5.80%  "
$ws.Range("D31").Value = "'1.180"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "'4.041"
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").Value = "'5.648"
$ws.Range("E33").Value = "  -16.93%  "
$ws.Range("D34").Value = "'0.08335"
$ws.Range("E34").Value = "  -5.42%  "
$ws.Range("D35").Value = "'1.672"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("D36").Value = "'12.31"
$ws.Range("E36").Value = "  -10.65%  "
$ws.Range("D37").Value = "'5.227"
$ws.Range("E37").Value = "  -6.72%  "
$ws.Range("D38").Value = "'0.06049"
$ws.Range("E38").Value = "  -7.49%  "
$ws.Range("D39").Value = "'0.02206"
$ws.Range("E39").Value = "  -7.81%  "
$ws.Range("D40").Value = "'1.205"
$ws.Range("E40").Value = "  -5.72%  "
$ws.Range("E41").Value = "  -6.66%  "
$ws.Range("D42").Value = "'8.194"
$ws.Range("E42").Value = "  -8.46%  "
$ws.Range("D43").Value = "'0.9985"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "'0.5843"
$ws.Range("E44").Value = "  -9.03%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("D46").Value = "'12.66"
$ws.Range("E46").Value = "  -9.22%  "
$ws.Range("D47").Value = "'0.5576"
$ws.Range("E47").Value = "  -8.98%  "
$ws.Range("D48").Value = "'122.09"
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("D49").Value = "'1.943"
$ws.Range("E49").Value = "  -8.99%  "
$ws.Range("D50").Value = "'0.06899"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("D51").Value = "'74.37"
$ws.Range("E51").Value = "  -6.68%  "
